# Apply odds corrections to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("H3").Value = 2.7
$ws.Range("V3").Value = 13
$ws.Range("AI3").Value = 26
# Row 4
$ws.Range("G4").Value = 2.2
$ws.Range("I4").Value = 3.1
$ws.Range("J4").Value = 1.07
$ws.Range("K4").Value = 9
$ws.Range("L4").Value = 1.36
$ws.Range("M4").Value = 3
$ws.Range("N4").Value = 2.1
$ws.Range("O4").Value = 1.7
$ws.Range("P4").Value = 1.44
$ws.Range("Q4").Value = 2.63
$ws.Range("T4").Value = 7
$ws.Range("Z4").Value = 9
$ws.Range("AB4").Value = 17
$ws.Range("AE4").Value = 8.5
# Row 5
$ws.Range("G5").Value = 4.5
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 1.57
$ws.Range("J5").Value = 1.01
$ws.Range("K5").Value = 23
$ws.Range("N5").Value = 1.44
$ws.Range("O5").Value = 2.7
$ws.Range("R5").Value = 1.53
$ws.Range("S5").Value = 2.38
$ws.Range("U5").Value = 29
$ws.Range("Y5").Value = 29
$ws.Range("AA5").Value = 10
$ws.Range("AB5").Value = 15
$ws.Range("AF5").Value = 10
$ws.Range("AH5").Value = 13
# Row 6
$ws.Range("G6").Value = 1.75
$ws.Range("H6").Value = 3.8
$ws.Range("I6").Value = 4.5
$ws.Range("AA6").Value = 7
$ws.Range("AB6").Value = 13
$ws.Range("AD6").Value = 151
$ws.Range("AE6").Value = 15
$ws.Range("AG6").Value = 15
$ws.Range("AH6").Value = 51
# Row 12
$ws.Range("G12").Value = 2.95
$ws.Range("I12").Value = 3.15
$ws.Range("T12").Value = 5.8
$ws.Range("U12").Value = 13
$ws.Range("W12").Value = 40
$ws.Range("Y12").Value = 60
$ws.Range("AE12").Value = 6.1
$ws.Range("AF12").Value = 14.5
$ws.Range("AH12").Value = 50
# Row 14
$ws.Range("N14").Value = 2.6
$ws.Range("O14").Value = 1.48
# Row 16
$ws.Range("J16").Value = 1.06
$ws.Range("K16").Value = 10
$ws.Range("L16").Value = 1.33
$ws.Range("M16").Value = 3.25
# Row 18
$ws.Range("G18").Value = 2.15
$ws.Range("I18").Value = 3.4
$ws.Range("K18").Value = 7.5
$ws.Range("T18").Value = 6.5
$ws.Range("X18").Value = 21
# Row 19
$ws.Range("G19").Value = 2.55
$ws.Range("I19").Value = 3.1
$ws.Range("U19").Value = 11
# Row 20
$ws.Range("N20").Value = 2.05
$ws.Range("O20").Value = 1.75
# Row 21
$ws.Range("G21").Value = 3.2
$ws.Range("I21").Value = 2.25
$ws.Range("J21").Value = 1.06
$ws.Range("K21").Value = 10
$ws.Range("N21").Value = 1.98
$ws.Range("O21").Value = 1.83
$ws.Range("U21").Value = 15
$ws.Range("AG21").Value = 9.5
$ws.Range("AI21").Value = 19
$ws.Range("AJ21").Value = 29
# Row 22
$ws.Range("N22").Value = 1.83
$ws.Range("O22").Value = 1.98
$ws.Range("AI22").Value = 21
# Row 23
$ws.Range("K23").Value = 8
# Row 24
$ws.Range("G24").Value = 2.75
$ws.Range("I24").Value = 2.75
$ws.Range("AB24").Value = 15
$ws.Range("AD24").Value = 351
$ws.Range("AF24").Value = 12
$ws.Range("AH24").Value = 26
$ws.Range("AI24").Value = 23
# Row 25
$ws.Range("N25").Value = 2.3
$ws.Range("O25").Value = 1.6
# Row 27
$ws.Range("J27").Value = 1.08
$ws.Range("K27").Value = 8
# Row 29
$ws.Range("L29").Value = 1.29
$ws.Range("M29").Value = 3.5
$ws.Range("N29").Value = 1.9
$ws.Range("O29").Value = 1.9
# Row 41
$ws.Range("J41").Value = 1.07
$ws.Range("K41").Value = 9
# Row 42
$ws.Range("N42").Value = 1.75
$ws.Range("O42").Value = 2.05

